$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 330, shifting existing rows 330-448 down to 331-449
$ws.Rows("330:330").Insert()

# Populate the newly inserted row 330 with the new data record
$ws.Range("A330").Value = 3
$ws.Range("B330").Value = "Femacal de La Calera"
$ws.Range("C330").Value = "Coquimbo"
$ws.Range("D330").Value = "2022-12-05"
$ws.Range("E330").Value = 5
$ws.Range("F330").Value = 100114013
$ws.Range("G330").Value = "Zanahoria"
$ws.Range("H330").Value = "Sin especificar"
$ws.Range("I330").Value = "Primera"
$ws.Range("J330").Value = 370
$ws.Range("K330").Value = 8000
$ws.Range("L330").Value = 8500
$ws.Range("M330").Value = 8216
$ws.Range("N330").Value = "$/saco 20 kilos"
$ws.Range("O330").Value = "Provincia de Quillota"
$ws.Range("P330").Value = 411
$ws.Range("Q330").Value = 20
$ws.Range("R330").Value = "Hortaliza"
